$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the reference-link lists in column D to reflect the refreshed
# department/team membership (strings changed, some rows trimmed).
$ws.Range("D2").Value = "UserDevelopment_Team_Denmark_Company,UserDevelopment_Team_Netherlands_Company,UserDevelopment_Team_Sweden_Company,UserDevelopment_Team_Germany_Company,Event,Europe_User_Development,Communication_and_PR"
$ws.Range("D3").Value = "UserTeam_Norway_Company,UserOperations_Team_Germany_Company,UserOperations_Team_Netherlands_Company,UserOperations_Team_Denmark_Company,UserOperations_Team_Sweden_Company,Fleet_Management_Team_Netherlands_Company,Fleet_Management_Team_Sweden_Company,Fleet_Management_Team_Denmark_Company,Business_Development,Business_Intelligence,Commercial_Product,Europe_Commercial_Operation,Europe_Business_Operation,Fleet_and_Business_Sales,Partner_Strategy,Sales_Planning"
$ws.Range("D4").Value = "UserRelationship_Team_Netherlands_Company,UserRelationship_Team_Denmark_Company,UserRelations_Team_Germany_Company"
$ws.Range("D5").Value = "Construction_Management,Design_Management,NIO_House_Operation,Europe_Space_Experience,PMO"
$ws.Range("D6").Value = "Power_Business_Operation_Team,Power_Management_Team_Germany_Company,Power_Management_Team_Netherlands_Company,Strategy and Business_Development Team,Market_Launch and Enabling Team,Power_Management_Team_Sweden_Company,Power_Management_Team_Denmark_Company,Europe_Power_Operation,Power_Operation,Power_market_launch_and_enabling_team"
$ws.Range("D7").Value = "Service_PMO,Service_Planning,Service_Operations_Team_Germany_Company,Service_Team_Norway_Company,Service_Operations_Team_Denmark_Company,Service_Operations_Team_Netherlands_Company,Service_Operations_Team_Sweden_Company,Europe_Service_Operation,Operation_Support,Parts_and_Logistics"
$ws.Range("D8").Value = "_Europe_Business_Enabling,Retall_Sales,Europe_Market_Planning"

# System_Development (row 10) and Training_Operation (row 11) no longer
# carry a department list.
$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()

# COE (row 12) keeps its original list text.
$ws.Range("D12").Value = "Europe_Business_HRBP_Department,Controlling_and_Planning_Department,Legal_EU_Department,NIO_Life_Supply_Chain_Department,Digital_Development_PMO_Team,Product_Marketing_Department,Europe_Product_Experience_Department,Purchasing_Governance_and_BP_Team"
